$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Set all new cell text values first, in the exact order the strings
#     were first introduced in the authored workbook, so the shared
#     string table comes out in the same order. ---
$ws.Range("A61").Value = "MARZO"
$ws.Range("B62").Value = "https://www.instagram.com/by.ornebad/"
$ws.Range("A62").Value = "by.ornebad"
$ws.Range("B63").Value = "https://www.instagram.com/jai._17/"
$ws.Range("A63").Value = "jai._17"
$ws.Range("B61").Value = "AGREGAR PORFOLIO DE FERNANDA GESE"
$ws.Range("E62").Value = "faltan textos"
$ws.Range("C62").Value = "fotografia"
$ws.Range("C63").Value = "pintura"

# --- Style-only fixes: H36 and H40 go from style "1" to style "2" ---
$ws.Range("H3").Copy()
$ws.Range("H36").PasteSpecial(-4122)
$ws.Range("H40").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 62 formatting (copy from the equivalent existing data row) ---
$ws.Range("D3").Copy()
$ws.Range("D62").PasteSpecial(-4122)
$ws.Range("E62").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E58").Copy()
$ws.Range("F62").PasteSpecial(-4122)
$ws.Range("G62").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Turn B62 into a real hyperlink pointing at the same address, then restore
# the usual hyperlink-cell formatting (style index 3 / "Hipervinculo")
$ws.Hyperlinks.Add($ws.Range("B62"), "https://www.instagram.com/by.ornebad/")
$ws.Range("B36").Copy()
$ws.Range("B62").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 63 formatting (copy from the equivalent existing data row) ---
$ws.Range("D3").Copy()
$ws.Range("D63").PasteSpecial(-4122)
$ws.Range("E63").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E58").Copy()
$ws.Range("F63").PasteSpecial(-4122)
$ws.Range("G63").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Hyperlinks.Add($ws.Range("B63"), "https://www.instagram.com/jai._17/")
$ws.Range("B36").Copy()
$ws.Range("B63").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the selected cell shown when the sheet is opened ---
$ws.Range("G40").Select() | Out-Null
